$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.668.60"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "3.885.34"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.10%  "
$c = $ws.Range("D5")
$c.Value = "'605.72"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$c = $ws.Range("D6")
$c.Value = "'170.09"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("D7").Value = "3.886.13"
$ws.Range("E8").Value = "  +0.17%  "
$c = $ws.Range("D9")
$c.Value = "'0.535"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.169"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.69%  "
$c = $ws.Range("D11")
$c.Value = "'6.39"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  +5.64%  "
$c = $ws.Range("D14")
$c.Value = "'38.30"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").Value = "4.543.49"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "3.902.37"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "69.616.95"
$ws.Range("E17").Value = "  +0.97%  "
$c = $ws.Range("D18")
$c.Value = "'18.66"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +9.03%  "
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("E20").Value = "  -0.68%  "
$c = $ws.Range("D21")
$c.Value = "'11.13"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.86%  "
$c = $ws.Range("D22")
$c.Value = "'490.58"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("E23").Value = "  +4.44%  "
$ws.Range("E24").Value = "  +3.55%  "
$c = $ws.Range("D25")
$c.Value = "'85.36"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +3.61%  "
$c = $ws.Range("D27")
$c.Value = "'12.33"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.10%  "
$c = $ws.Range("D28")
$c.Value = "'10.15"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("E30").Value = "  +1.19%  "
$c = $ws.Range("D31")
$c.Value = "'2.42"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").Value = "4.038.85"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  -0.46%  "
$c = $ws.Range("D34")
$c.Value = "'31.92"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").Value = "3.851.53"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").Value = "  +0.73%  "
$c = $ws.Range("D37")
$c.Value = "'6.13"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.51%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +1.14%  "
$c = $ws.Range("D40")
$c.Value = "'3.32"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +12.05%  "
$c = $ws.Range("D41")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "
$c = $ws.Range("D42")
$c.Value = "'0.328"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("E43").Value = "  +6.49%  "
$c = $ws.Range("D44")
$c.Value = "'436.99"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.32%  "
$c = $ws.Range("D46")
$c.Value = "'8.71"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.09%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +3.33%  "
$c = $ws.Range("D49")
$c.Value = "'0.000275"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +21.30%  "
$c = $ws.Range("D50")
$c.Value = "'143.95"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$c = $ws.Range("D51")
$c.Value = "'40.23"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.03%  "
